# Generate Report for Handback
#
# Semantic changes (per the OOXML diff):
#  1. Status "In Translation" -> "Handed back: in sync with en-US"
#     (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3 all share this string)
#  2. zh-cn / de-de sheets: populate the "Latest Target File" (I),
#     "Latest Handback File" (J) and "Latest Handback DateTime" (K)
#     columns for the two data rows, now that handback has happened.
#     "Latest Target File" becomes a hyperlink to the source .md file
#     (same target/display text already used by column A's hyperlink).
#  3. Column width bump on the columns whose text just got much longer
#     (Status column + the two newly-populated file-name columns).

$wb = $excel.ActiveWorkbook

$mdAddr1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d221d5eaa2c6eaadee715e95938442594440caa6/e2e/0036e3ec-ce58-4674-92eb-70005cb62c1b.md"
$mdDisp1 = "0036e3ec-ce58-4674-92eb-70005cb62c1b.md"
$mdAddr2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d221d5eaa2c6eaadee715e95938442594440caa6/e2e/8fae9750-eb17-4b65-bd2c-a489042702eb.md"
$mdDisp2 = "8fae9750-eb17-4b65-bd2c-a489042702eb.md"

$newStatus = "Handed back: in sync with en-US"

# --- 1. Overview sheet: Status columns (E/F) for both rows ---------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

# --- 2. zh-cn sheet --------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdAddr1, [Type]::Missing, [Type]::Missing, $mdDisp1)
$wsZh.Range("J2").Value = "0036e3ec-ce58-4674-92eb-70005cb62c1b.98f095102b573fc4c9176676c41fb5bab78a0601.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-01 12:26:11"

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdAddr2, [Type]::Missing, [Type]::Missing, $mdDisp2)
$wsZh.Range("J3").Value = "8fae9750-eb17-4b65-bd2c-a489042702eb.142369e34c7abe1ddcaf83f2b49806a6d61b60ff.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-01 12:26:11"

$wsZh.Columns.Item(3).ColumnWidth = 29.15
$wsZh.Columns.Item(9).ColumnWidth = 39.1666666666667
$wsZh.Columns.Item(10).ColumnWidth = 39.1666666666667

# --- 3. de-de sheet ---------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdAddr1, [Type]::Missing, [Type]::Missing, $mdDisp1)
$wsDe.Range("J2").Value = "0036e3ec-ce58-4674-92eb-70005cb62c1b.98f095102b573fc4c9176676c41fb5bab78a0601.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-01 12:26:21"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdAddr2, [Type]::Missing, [Type]::Missing, $mdDisp2)
$wsDe.Range("J3").Value = "8fae9750-eb17-4b65-bd2c-a489042702eb.142369e34c7abe1ddcaf83f2b49806a6d61b60ff.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-01 12:26:21"

$wsDe.Columns.Item(3).ColumnWidth = 29.15
$wsDe.Columns.Item(9).ColumnWidth = 39.1666666666667
$wsDe.Columns.Item(10).ColumnWidth = 39.1666666666667
